$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.231350064277649
$ws.Range("B1").Value = 2.511218309402466
$ws.Range("C1").Value = 4.514007091522217
$ws.Range("D1").Value = 2.51203465461731
$ws.Range("E1").Value = 1.073765993118286
